$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "#51 FANCY SEA MAXI DRESS - ANDRE"
$ws.Range("B2").Value = "#88 WESTERN WAVE MIDI SKIRT"

$ws.Range("A3").Value = "#138 Y #140  SLIP SKIRT"
$ws.Range("B3").Value = "#121 COLLAR DXF"

$ws.Range("A4").Value = "#88 WESTERN WAVE MIDI SKIRT"
$ws.Range("B4").Value = "#24 COLLAR DXF"

$ws.Range("A5").Value = "#100 BOAT LINES TOP-LUZKA"
$ws.Range("B5").Value = "#79 COLLAR -NECKLACE DXF"

$ws.Range("A6").Value = "#121 COLLAR DXF"
$ws.Range("B6").Value = ""

$ws.Range("A7").Value = "#24 COLLAR DXF"
$ws.Range("B7").Value = ""

$ws.Range("A8").Value = "#79 COLLAR -NECKLACE DXF"
$ws.Range("B8").Value = ""
